# Student Onboarding - Work Experience
# 1) Rename the placeholder "Sheet8" tab to "WorkExperience"
# 2) Refresh the demo student's name / email on the RegistrationForm sheet
#    and hyperlink the new email address
# 3) Populate the new WorkExperience sheet with its header row + one sample
#    row of data, matching the other "…Details" sheets already in the file

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename Sheet8 -> WorkExperience
# ---------------------------------------------------------------------
$work = $wb.Worksheets.Item("Sheet8")
$work.Name = "WorkExperience"

# ---------------------------------------------------------------------
# 2) RegistrationForm: swap demo name + bump the demo e-mail, then
#    re-hyperlink the e-mail cell to the new address
# ---------------------------------------------------------------------
$reg = $wb.Worksheets.Item("RegistrationForm")
$reg.Activate()

$reg.Range("A2").Value = "Prem"

$newEmail = "priya.t+studentdemo68@icanio.com"
$reg.Range("C2").Value = $newEmail
$reg.Hyperlinks.Add($reg.Range("C2"), "mailto:" + $newEmail, "", "", $newEmail)

# ---------------------------------------------------------------------
# 3) WorkExperience: header row + first data row
# ---------------------------------------------------------------------
$work.Activate()

$headers = @("Jobtitle","Companyname","Startmonth","Startyear","Noticeperiod","Function","Role","Industry","Description","Skill1","Skill2","Skill3","Awardname","AwardReceivedname","AwardDescription","bestperformancecertificate")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $work.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$work.Range("A2").Value = "FrontEnd Developer"
$work.Range("B2").Value = "Icanio Technologies 5.0 "
$work.Range("C2").Value = "August"
$work.Range("D2").Value = 2022
$work.Range("E2").Value = 10
$work.Range("F2").Value = "Frontend Developer "
$work.Range("G2").Value = "Junior Web Developer "
$work.Range("H2").Value = "Software Industry "
$work.Range("I2").Value = "Need a career growth. "
$work.Range("J2").Value = "Nodejs "
$work.Range("K2").Value = "Javascript "
$work.Range("L2").Value = "Angularjs "
$work.Range("M2").Value = "Best Performance"
$work.Range("N2").Value = 45530
$work.Range("O2").Value = "For my work effect received a Best performer award."
$work.Range("P2").Value = "C:\Users\ICANIO-10090\Pictures\Bestperformance certificate.jpg"

# Wrap the longer free-text fields, matching the look of the row-1 wrapped
# cells used elsewhere in the workbook, and size row 2 to fit.
$work.Range("B2").WrapText = $true
$work.Range("J2").WrapText = $true
$work.Range("K2").WrapText = $true
$work.Range("L2").WrapText = $true
$work.Rows.Item(2).RowHeight = 35

$work.Range("D4").Select()

# ---------------------------------------------------------------------
# Leave RegistrationForm as the active sheet/selection, matching the
# workbook's last-saved cursor position
# ---------------------------------------------------------------------
$reg.Activate()
$reg.Range("N2").Select()
